# 2.10.7 Remove requiredVL parameter
$wb = $excel.ActiveWorkbook

# Delete the "requiredvl" row from the "Model parameters" sheet (row 76)
$wsModelParams = $wb.Worksheets.Item("Model parameters")
$wsModelParams.Rows.Item(76).Delete()

# Delete the "requiredvl" row from the "Data constants" sheet (row 25)
$wsDataConstants = $wb.Worksheets.Item("Data constants")
$wsDataConstants.Rows.Item(25).Delete()

# Set the active sheet / view state to match the resulting workbook state
$wsModelParams.Activate()
$wsModelParams.Application.ActiveWindow.ScrollRow = 2
$wsModelParams.Range("L28").Select()
